# food-details 코드 오류 해결
# Populate the "Lunch" worksheet with the corrected food-details rows
# (fixes the "Ingrdient Code" -> "Ingredient Code" header typo and adds
# the missing ingredient breakdown rows for this sheet's dishes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lunch")

# Header row
$ws.Range("A1").Value = "Food Code"
$ws.Range("B1").Value = "Food Name"
$ws.Range("C1").Value = "Ingredient Code"
$ws.Range("D1").Value = "Ingrdient"
$ws.Range("E1").Value = "1 person (g)"

# Caramelized Pork with Eggs
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Caramelized Pork with Eggs"
$ws.Range("C2").Value = "H0002"
$ws.Range("D2").Value = "Pork meat, lean, raw"
$ws.Range("E2").Value = 100

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Caramelized Pork with Eggs"
$ws.Range("C3").Value = "B0008"
$ws.Range("D3").Value = "Garlic, fresh, raw"
$ws.Range("E3").Value = 6.66

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Caramelized Pork with Eggs"
$ws.Range("C4").Value = "N0001"
$ws.Range("D4").Value = "Salt, table"
$ws.Range("E4").Value = 3.33

# Num Banchok
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "Num Banchok"
$ws.Range("C5").Value = "N0005"
$ws.Range("D5").Value = "Noodle, rice flour, wet"
$ws.Range("E5").Value = 61

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Num Banchok"
$ws.Range("C6").Value = "J0012"
$ws.Range("D6").Value = "Fish, Mystus wolffi, raw"
$ws.Range("E6").Value = 80

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Num Banchok"
$ws.Range("C7").Value = "G0001"
$ws.Range("D7").Value = "Sugar, granulated"
$ws.Range("E7").Value = 5

# Match the bold/centered/bordered header style already used by the
# other sheets' header rows (reuses the same style record instead of
# creating a near-duplicate one).
$headerSrc = $wb.Worksheets.Item("Breakfast").Range("A1:E1")
$headerSrc.Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)
